# Refresh the "cryptos" price/volume table to the latest scrape.
# Cells whose new Price text is a plain decimal (e.g. "212.89") would be
# auto-converted to a Number by COM's normal type inference, so those are
# written with a leading apostrophe (forces text entry, like typing it into
# Excel) and then have their style reset back to "Normal" so no stray
# NumberFormat/quote-prefix styling is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.643.75"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "1.636.67"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.89"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").Value = "1.865.92"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "1.647.93"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "26.636.12"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "'63.40"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "'218.13"
$ws.Range("E19").Value = "  +7.49%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "'9.51"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("E23").Value = "  +3.23%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'149.18"
$ws.Range("E25").Value = "  +4.75%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "'6.89"
$ws.Range("E28").Value = "  +5.07%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0174"
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.182.09"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").Value = "'0.810"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.508"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'5.42"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "'0.794"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.773.51"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "'92.67"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "'54.89"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").Value = "'7.65"
$ws.Range("E49").Value = "  +4.89%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +0.11%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
